$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 16 (everything from old row 16 downward shifts to 17+)
$ws.Rows("16:16").Insert()

$ws.Range("A16").Value = 1
$ws.Range("B16").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C16").Value = "Arica y Parinacota"
$ws.Range("D16").Value = 44453
$ws.Range("E16").Value = 15
$ws.Range("F16").Value = 100112028
$ws.Range("G16").Value = "Sandia"
$ws.Range("H16").Value = "Sin especificar"
$ws.Range("I16").Value = "Tercera"
$ws.Range("J16").Value = 700
$ws.Range("K16").Value = 800
$ws.Range("L16").Value = 850
$ws.Range("M16").Value = 825
$ws.Range("N16").Value = "`$/kilo (volumen en unidades)"
$ws.Range("O16").Value = "Perú"
$ws.Range("P16").Value = 825
$ws.Range("Q16").Value = 1
$ws.Range("R16").Value = "Hortaliza"

# Insert another new row at position 18 (everything from old row 17 - now at 17 - downward from 18 shifts to 19+)
$ws.Rows("18:18").Insert()

$ws.Range("A18").Value = 1
$ws.Range("B18").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C18").Value = "Arica y Parinacota"
$ws.Range("D18").Value = 44523
$ws.Range("E18").Value = 15
$ws.Range("F18").Value = 100112028
$ws.Range("G18").Value = "Sandia"
$ws.Range("H18").Value = "Sin especificar"
$ws.Range("I18").Value = "Segunda"
$ws.Range("J18").Value = 1000
$ws.Range("K18").Value = 550
$ws.Range("L18").Value = 580
$ws.Range("M18").Value = 565
$ws.Range("N18").Value = "`$/kilo (volumen en unidades)"
$ws.Range("O18").Value = "Perú"
$ws.Range("P18").Value = 565
$ws.Range("Q18").Value = 1
$ws.Range("R18").Value = "Hortaliza"

# Ensure the date cells keep the date-formatted style (style index 2) used by the rest of column D
$ws.Range("D16").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("D18").NumberFormat = "YYYY-MM-DD HH:MM:SS"
